# Updates TPM-derived NATMI metrics for sheet1 (C3-Cr2 LR-pairs)
# per commit "update scripts wuth new tpm": refreshed ligand/receptor
# expression + specificity + edge-weight values recalculated from new TPM input.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.555934
$ws.Range("H2").Value = 1.667802
$ws.Range("I2").Value = 0.005745252779589096
$ws.Range("J2").Value = 0.005745252779589094
$ws.Range("M2").Value = 1.028288333333333
$ws.Range("N2").Value = 3.084865
$ws.Range("O2").Value = 0.5720958559345812
$ws.Range("P2").Value = 0.5720958559345812
$ws.Range("Q2").Value = 0.5716604463033333
$ws.Range("R2").Value = 5.144944016729999
$ws.Range("S2").Value = 0.003286835306499556
$ws.Range("T2").Value = 0.003286835306499555

# Row 3
$ws.Range("G3").Value = 0.555934
$ws.Range("H3").Value = 1.667802
$ws.Range("I3").Value = 0.005745252779589096
$ws.Range("J3").Value = 0.005745252779589094
$ws.Range("O3").Value = 0.2417269186310566
$ws.Range("P3").Value = 0.2417269186310566
$ws.Range("Q3").Value = 0.2415429455653333
$ws.Range("R3").Value = 2.173886510088
$ws.Range("S3").Value = 0.001388782251166585
$ws.Range("T3").Value = 0.001388782251166585

# Row 4
$ws.Range("G4").Value = 0.555934
$ws.Range("H4").Value = 1.667802
$ws.Range("I4").Value = 0.005745252779589096
$ws.Range("J4").Value = 0.005745252779589094
$ws.Range("O4").Value = 0.03814534912077908
$ws.Range("P4").Value = 0.03814534912077908
$ws.Range("Q4").Value = 0.03811631753066667
$ws.Range("R4").Value = 0.343046857776
$ws.Range("S4").Value = 0.0002191546730645525
$ws.Range("T4").Value = 0.0002191546730645524

# Row 5
$ws.Range("G5").Value = 0.555934
$ws.Range("H5").Value = 1.667802
$ws.Range("I5").Value = 0.005745252779589096
$ws.Range("J5").Value = 0.005745252779589094
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.2660733333333333
$ws.Range("N5").Value = 0.7982199999999999
$ws.Range("O5").Value = 0.1480318763135831
$ws.Range("P5").Value = 0.1480318763135831
$ws.Range("Q5").Value = 0.1479192124933333
$ws.Range("R5").Value = 1.33127291244
$ws.Range("S5").Value = 0.0008504805488584024
$ws.Range("T5").Value = 0.0008504805488584021

# Row 6
$ws.Range("I6").Value = 0.823525905561055
$ws.Range("J6").Value = 0.823525905561055
$ws.Range("M6").Value = 1.028288333333333
$ws.Range("N6").Value = 3.084865
$ws.Range("O6").Value = 0.5720958559345812
$ws.Range("P6").Value = 0.5720958559345812
$ws.Range("Q6").Value = 81.94194490238944
$ws.Range("R6").Value = 737.4775041215049
$ws.Range("S6").Value = 0.4711357578262529
$ws.Range("T6").Value = 0.4711357578262529

# Row 7
$ws.Range("I7").Value = 0.823525905561055
$ws.Range("J7").Value = 0.823525905561055
$ws.Range("O7").Value = 0.2417269186310566
$ws.Range("P7").Value = 0.2417269186310566
$ws.Range("S7").Value = 0.1990683795641243
$ws.Range("T7").Value = 0.1990683795641243

# Row 8
$ws.Range("I8").Value = 0.823525905561055
$ws.Range("J8").Value = 0.823525905561055
$ws.Range("O8").Value = 0.03814534912077908
$ws.Range("P8").Value = 0.03814534912077908
$ws.Range("S8").Value = 0.03141368317763218
$ws.Range("T8").Value = 0.03141368317763218

# Row 9
$ws.Range("I9").Value = 0.823525905561055
$ws.Range("J9").Value = 0.823525905561055
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.2660733333333333
$ws.Range("N9").Value = 0.7982199999999999
$ws.Range("O9").Value = 0.1480318763135831
$ws.Range("P9").Value = 0.1480318763135831
$ws.Range("Q9").Value = 21.20277524623778
$ws.Range("R9").Value = 190.82497721614
$ws.Range("S9").Value = 0.1219080849930456
$ws.Range("T9").Value = 0.1219080849930456

# Row 10
$ws.Range("G10").Value = 0.3446996666666666
$ws.Range("H10").Value = 1.034099
$ws.Range("I10").Value = 0.003562269474506148
$ws.Range("J10").Value = 0.003562269474506148
$ws.Range("M10").Value = 1.028288333333333
$ws.Range("N10").Value = 3.084865
$ws.Range("O10").Value = 0.5720958559345812
$ws.Range("P10").Value = 0.5720958559345812
$ws.Range("Q10").Value = 0.3544506457372221
$ws.Range("R10").Value = 3.190055811634999
$ws.Range("S10").Value = 0.002037959604087226
$ws.Range("T10").Value = 0.002037959604087225

# Row 11
$ws.Range("G11").Value = 0.3446996666666666
$ws.Range("H11").Value = 1.034099
$ws.Range("I11").Value = 0.003562269474506148
$ws.Range("J11").Value = 0.003562269474506148
$ws.Range("O11").Value = 0.2417269186310566
$ws.Range("P11").Value = 0.2417269186310566
$ws.Range("Q11").Value = 0.1497655707728889
$ws.Range("R11").Value = 1.347890136956
$ws.Range("S11").Value = 0.0008610964234058443
$ws.Range("T11").Value = 0.0008610964234058443

# Row 12
$ws.Range("G12").Value = 0.3446996666666666
$ws.Range("H12").Value = 1.034099
$ws.Range("I12").Value = 0.003562269474506148
$ws.Range("J12").Value = 0.003562269474506148
$ws.Range("O12").Value = 0.03814534912077908
$ws.Range("P12").Value = 0.03814534912077908
$ws.Range("Q12").Value = 0.02363352834577778
$ws.Range("R12").Value = 0.212701755112
$ws.Range("S12").Value = 0.0001358840127673313
$ws.Range("T12").Value = 0.0001358840127673313

# Row 13
$ws.Range("G13").Value = 0.3446996666666666
$ws.Range("H13").Value = 1.034099
$ws.Range("I13").Value = 0.003562269474506148
$ws.Range("J13").Value = 0.003562269474506148
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.2660733333333333
$ws.Range("N13").Value = 0.7982199999999999
$ws.Range("O13").Value = 0.1480318763135831
$ws.Range("P13").Value = 0.1480318763135831
$ws.Range("Q13").Value = 0.09171538930888888
$ws.Range("R13").Value = 0.8254385037799998
$ws.Range("S13").Value = 0.0005273294342457468
$ws.Range("T13").Value = 0.0005273294342457467

# Row 14
$ws.Range("G14").Value = 16.17571666666667
$ws.Range("H14").Value = 48.52715
$ws.Range("I14").Value = 0.1671665721848498
$ws.Range("J14").Value = 0.1671665721848498
$ws.Range("M14").Value = 1.028288333333333
$ws.Range("N14").Value = 3.084865
$ws.Range("O14").Value = 0.5720958559345812
$ws.Range("P14").Value = 0.5720958559345812
$ws.Range("Q14").Value = 16.63330073163889
$ws.Range("R14").Value = 149.69970658475
$ws.Range("S14").Value = 0.09563530319774163
$ws.Range("T14").Value = 0.09563530319774162

# Row 15
$ws.Range("G15").Value = 16.17571666666667
$ws.Range("H15").Value = 48.52715
$ws.Range("I15").Value = 0.1671665721848498
$ws.Range("J15").Value = 0.1671665721848498
$ws.Range("O15").Value = 0.2417269186310566
$ws.Range("P15").Value = 0.2417269186310566
$ws.Range("Q15").Value = 7.028046944955555
$ws.Range("R15").Value = 63.2524225046
$ws.Range("S15").Value = 0.04040866039235985
$ws.Range("T15").Value = 0.04040866039235985

# Row 16
$ws.Range("G16").Value = 16.17571666666667
$ws.Range("H16").Value = 48.52715
$ws.Range("I16").Value = 0.1671665721848498
$ws.Range("J16").Value = 0.1671665721848498
$ws.Range("O16").Value = 0.03814534912077908
$ws.Range("P16").Value = 0.03814534912077908
$ws.Range("Q16").Value = 1.109050269911111
$ws.Range("R16").Value = 9.981452429200001
$ws.Range("S16").Value = 0.006376627257315015
$ws.Range("T16").Value = 0.006376627257315014

# Row 17
$ws.Range("G17").Value = 16.17571666666667
$ws.Range("H17").Value = 48.52715
$ws.Range("I17").Value = 0.1671665721848498
$ws.Range("J17").Value = 0.1671665721848498
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.2660733333333333
$ws.Range("N17").Value = 0.7982199999999999
$ws.Range("O17").Value = 0.1480318763135831
$ws.Range("P17").Value = 0.1480318763135831
$ws.Range("Q17").Value = 4.303926852555556
$ws.Range("R17").Value = 38.73534167299999
$ws.Range("S17").Value = 0.02474598133743335
$ws.Range("T17").Value = 0.02474598133743335
